$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = 0.8396972991405803
$ws.Range("C6").Value = 0.0312041741843776
$ws.Range("D6").Value = 0.7987149222933083
$ws.Range("E6").Value = 0.8107406403421583
$ws.Range("F6").Value = 0.8572729988052569
$ws.Range("G6").Value = 0.8475678443420379
$ws.Range("H6").Value = 0.8841900899201408
$ws.Range("I6").Value = 0.8233575453282475
$ws.Range("J6").Value = 0.0285476826768353
$ws.Range("K6").Value = 0.7933452821461611
$ws.Range("L6").Value = 0.7896813759716984
$ws.Range("M6").Value = 0.8460604917311199
$ws.Range("N6").Value = 0.8253475616378843
$ws.Range("O6").Value = 0.8623530151543735
$ws.Range("P6").Value = 0.7324787718191119
$ws.Range("Q6").Value = 0.03013145022884374
$ws.Range("R6").Value = 0.702901543596162
$ws.Range("S6").Value = 0.7243672115139208
$ws.Range("T6").Value = 0.7327680454802696
$ws.Range("U6").Value = 0.7131031691811918
$ws.Range("V6").Value = 0.7892538893240154
$ws.Range("W6").Value = 0.807798954614249
$ws.Range("X6").Value = 0.01904189432235215
$ws.Range("Y6").Value = 0.7776350127413958
$ws.Range("Z6").Value = 0.8053209309139101
$ws.Range("AA6").Value = 0.8193198048885656
$ws.Range("AB6").Value = 0.8019098320331717
$ws.Range("AC6").Value = 0.834809192494202
$ws.Range("AD6").Value = 0.8017485226484183
$ws.Range("AE6").Value = 0.02187083570380738
$ws.Range("AF6").Value = 0.7861734635827251
$ws.Range("AG6").Value = 0.7917073965461062
$ws.Range("AH6").Value = 0.8086021505376344
$ws.Range("AI6").Value = 0.7809199386553173
$ws.Range("AJ6").Value = 0.841339663920309
$ws.Range("AK6").Value = 0.825416323160782
$ws.Range("AL6").Value = 0.04874581339422589
$ws.Range("AM6").Value = 0.7876899696048633
$ws.Range("AN6").Value = 0.7581132322593422
$ws.Range("AO6").Value = 0.8578255675029868
$ws.Range("AP6").Value = 0.8282556750298686
$ws.Range("AQ6").Value = 0.8951971714068487
$ws.Range("B7").Value = 0.8487159913715232
$ws.Range("C7").Value = 0.04655246907341332
$ws.Range("D7").Value = 0.8484393869529895
$ws.Range("F7").Value = 0.8749476738441084
$ws.Range("I7").Value = 0.8606353889948695
$ws.Range("J7").Value = 0.03018071514868441
$ws.Range("M7").Value = 0.863684452275534
$ws.Range("N7").Value = 0.8599590373783923
$ws.Range("P7").Value = 0.8396140992435329
$ws.Range("Q7").Value = 0.04151149063805919
$ws.Range("R7").Value = 0.827851481017939
$ws.Range("S7").Value = 0.8042205054064636
$ws.Range("T7").Value = 0.8435288050654922
$ws.Range("V7").Value = 0.9172469704727769
$ws.Range("W7").Value = 0.8497538216940527
$ws.Range("X7").Value = 0.03232703380369924
$ws.Range("Y7").Value = 0.845905867182463
$ws.Range("AB7").Value = 0.8537326251580331
$ws.Range("AC7").Value = 0.8723050396437493
$ws.Range("AD7").Value = 0.8647471614363464
$ws.Range("AE7").Value = 0.03459058903636719
$ws.Range("AF7").Value = 0.8799547969760736
$ws.Range("AG7").Value = 0.8021438796231567
$ws.Range("AH7").Value = 0.8736047107014848
$ws.Range("AJ7").Value = 0.9065001212709193
$ws.Range("AK7").Value = 0.8420709957133198
$ws.Range("AL7").Value = 0.02553522393627148
$ws.Range("AM7").Value = 0.8370406189555125
$ws.Range("AN7").Value = 0.8078975110134604
$ws.Range("AQ7").Value = 0.8865207373271887

Write-Host "Updated 72 cells in rows 6 and 7"
